$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: C7 0.6 -> 0.8
$ws.Range("C7").Value = 0.8

# Row 9: C9 0.5 -> 0.8
$ws.Range("C9").Value = 0.8

# Row 11: C11 0.5 -> 0.8
$ws.Range("C11").Value = 0.8

# Row 13: C13 0.6 -> 0.9, D13 text updated
$ws.Range("C13").Value = 0.9
$ws.Range("D13").Value = "Silhouette-score et sample silhouette"

# Row 14: C14 0.5 -> 0.8, D14 new text
$ws.Range("C14").Value = 0.8
$ws.Range("D14").Value = "Ari_Score sur 12 mois"

# Row 15: C15 0.5 -> 0.8
$ws.Range("C15").Value = 0.8

# Row 18: C18 0.5 -> 0.7
$ws.Range("C18").Value = 0.7

# Update selection
$ws.Range("B19:C19").Select()
